# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.837.44'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '2.206.98'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '287.94'
$ws.Range('E5').Value = '  -2.01%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '87.06'
$ws.Range('E6').Value = '  +3.16%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.513'
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('E8').Value = '  -0.03%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.468'
$ws.Range('E9').Value = '  +0.01%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '30.30'
$ws.Range('E10').Value = '  +1.25%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0775'
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('E12').Value = '  +2.19%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '6.41'
$ws.Range('E13').Value = '  +1.43%  '
$ws.Range('D14').Value = '2.548.17'
$ws.Range('E14').Value = '  -0.83%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '13.87'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').Value = '2.188.01'
$ws.Range('E16').Value = '  -1.55%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.724'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').Value = '39.775.82'
$ws.Range('E18').Value = '  +0.14%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.51'
$ws.Range('E19').Value = '  +9.20%  '
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('E21').Value = '  -0.08%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '65.30'
$ws.Range('E22').Value = '  -0.09%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '235.03'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('E26').Value = '  -1.26%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '22.37'
$ws.Range('E27').Value = '  -2.34%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.10'
$ws.Range('E28').Value = '  -0.98%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.18'
$ws.Range('E29').Value = '  -0.44%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '152.66'
$ws.Range('E30').Value = '  +1.33%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '31.54'
$ws.Range('E31').Value = '  -2.83%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.998'
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('E33').Value = '  +1.68%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.0714'
$ws.Range('E34').Value = '  +1.35%  '
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('E37').Value = '  -0.10%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '15.64'
$ws.Range('E38').Value = '  -2.81%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0979'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('E40').Value = '  +1.86%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.089.43'
$ws.Range('E41').Value = '  +7.33%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.83'
$ws.Range('E42').Value = '  +3.31%  '
$ws.Range('E43').Value = '  -0.53%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '9.92'
$ws.Range('E44').Value = '  +5.52%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0266'
$ws.Range('E45').Value = '  -0.45%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '17.37'
$ws.Range('E46').Value = '  +6.82%  '
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('D48').Value = '2.422.74'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '88.21'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '68.74'
$ws.Range('E50').Value = '  -2.94%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.44'
$ws.Range('E51').Value = '  +0.27%  '
